# Auto-generated PowerShell COM-interop script
# Applies numeric cell updates across the ALC..WVR sheets per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 295.33334
$ws.Range("I12").Value = 295.5
$ws.Range("J12").Value = 295
$ws.Range("K12").Value = 295.5
$ws.Range("L12").Value = 295
$ws.Range("M12").Value = -125.5
$ws.Range("N12").Value = -635
$ws.Range("H42").Value = 3412.25
$ws.Range("I42").Value = 643.75
$ws.Range("J42").Value = 4796.5
$ws.Range("K42").Value = 1931.25
$ws.Range("L42").Value = 14389.5
$ws.Range("M42").Value = -1701.25
$ws.Range("N42").Value = -14849.5
$ws.Range("H53").Value = 372.85715
$ws.Range("I53").Value = 277.625
$ws.Range("J53").Value = 499.83334
$ws.Range("K53").Value = 277.625
$ws.Range("L53").Value = 499.83334
$ws.Range("M53").Value = 359.375
$ws.Range("N53").Value = -1773.83334
$ws.Range("H55").Value = 69
$ws.Range("J55").Value = 69
$ws.Range("L55").Value = 69
$ws.Range("N55").Value = -497
$ws.Range("H98").Value = 1598.8
$ws.Range("I98").Value = 1598.8
$ws.Range("K98").Value = 1598.8
$ws.Range("M98").Value = -100.8
$ws.Range("H113").Value = 4480
$ws.Range("H116").Value = 5753.6
$ws.Range("I116").Value = 4980
$ws.Range("J116").Value = 5947
$ws.Range("K116").Value = 4980
$ws.Range("L116").Value = 5947
$ws.Range("M116").Value = -1538
$ws.Range("N116").Value = -12831
$ws.Range("H122").Value = 1598.8
$ws.Range("I122").Value = 1598.8
$ws.Range("K122").Value = 4796.4
$ws.Range("M122").Value = -2346.4
$ws.Range("H125").Value = 21248.25
$ws.Range("I125").Value = 18993
$ws.Range("J125").Value = 22000
$ws.Range("K125").Value = 170937
$ws.Range("L125").Value = 198000
$ws.Range("M125").Value = -168477
$ws.Range("N125").Value = -202920
$ws.Range("H132").Value = 3253
$ws.Range("I132").Value = 3253
$ws.Range("K132").Value = 9759
$ws.Range("M132").Value = -7229
$ws.Range("H138").Value = 2746.3076
$ws.Range("J138").Value = 4249.5
$ws.Range("L138").Value = 12748.5
$ws.Range("N138").Value = -23028.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4500
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H32").Value = 3211179
$ws.Range("I32").Value = 3046447.8
$ws.Range("K32").Value = 3046447.8
$ws.Range("M32").Value = -3046160.8
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H110").Value = 1082.5
$ws.Range("I110").Value = 1082.5
$ws.Range("K110").Value = 1082.5
$ws.Range("M110").Value = 962.5
$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H122").Value = 3333
$ws.Range("I122").Value = 3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9999
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("M122").Value = -7549
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4500
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H20").Value = 6400
$ws.Range("I20").Value = 4000
$ws.Range("K20").Value = 4000
$ws.Range("M20").Value = -3753
$ws.Range("H22").Value = 142.54546
$ws.Range("I22").Value = 136.6
$ws.Range("K22").Value = 136.6
$ws.Range("M22").Value = 36.40000000000001
$ws.Range("H86").Value = 1984.125
$ws.Range("I86").Value = 1828.8334
$ws.Range("J86").Value = 2450
$ws.Range("K86").Value = 1828.8334
$ws.Range("L86").Value = 2450
$ws.Range("M86").Value = -705.8334
$ws.Range("N86").Value = -4696
$ws.Range("H89").Value = 1984.125
$ws.Range("I89").Value = 1828.8334
$ws.Range("J89").Value = 2450
$ws.Range("K89").Value = 9144.166999999999
$ws.Range("L89").Value = 12250
$ws.Range("M89").Value = -3528.166999999999
$ws.Range("N89").Value = -23482
$ws.Range("H99").Value = 998
$ws.Range("I99").Value = 998
$ws.Range("K99").Value = 998
$ws.Range("M99").Value = 500
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1981.3334
$ws.Range("I31").Value = 1981.3334
$ws.Range("K31").Value = 1981.3334
$ws.Range("M31").Value = -1686.3334
$ws.Range("H34").Value = 1981.3334
$ws.Range("I34").Value = 1981.3334
$ws.Range("K34").Value = 1981.3334
$ws.Range("M34").Value = -1779.3334
$ws.Range("H122").Value = 1300
$ws.Range("J122").Value = 1300
$ws.Range("L122").Value = 3900
$ws.Range("N122").Value = -8800
$ws.Range("H132").Value = 5493.75
$ws.Range("I132").Value = 5992.6665
$ws.Range("K132").Value = 17977.9995
$ws.Range("M132").Value = -15447.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 96.066666
$ws.Range("J12").Value = 65.25
$ws.Range("L12").Value = 195.75
$ws.Range("N12").Value = -541.75
$ws.Range("H14").Value = 341.85715
$ws.Range("I14").Value = 341.85715
$ws.Range("K14").Value = 1025.57145
$ws.Range("M14").Value = -852.5714499999999
$ws.Range("H68").Value = 378.66666
$ws.Range("J68").Value = 413
$ws.Range("L68").Value = 1239
$ws.Range("N68").Value = -2861
$ws.Range("H71").Value = 378.66666
$ws.Range("J71").Value = 413
$ws.Range("L71").Value = 3717
$ws.Range("N71").Value = -11829
$ws.Range("H92").Value = 150
$ws.Range("H95").Value = 8946.5
$ws.Range("J95").Value = 8946.5
$ws.Range("L95").Value = 26839.5
$ws.Range("N95").Value = -30957.5
$ws.Range("H107").Value = 500
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340
$ws.Range("H118").Value = 1897.75
$ws.Range("I118").Value = 1897.75
$ws.Range("K118").Value = 5693.25
$ws.Range("M118").Value = -4450.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.916664
$ws.Range("I2").Value = 65.71429000000001
$ws.Range("K2").Value = 65.71429000000001
$ws.Range("M2").Value = 47.28570999999999
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H102").Value = 2996
$ws.Range("I102").Value = 2996
$ws.Range("K102").Value = 2996
$ws.Range("M102").Value = -1374
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 3774.75
$ws.Range("I122").Value = 3774.75
$ws.Range("K122").Value = 11324.25
$ws.Range("M122").Value = -8874.25
$ws.Range("H126").Value = 8325
$ws.Range("J126").Value = 8500
$ws.Range("L126").Value = 25500
$ws.Range("N126").Value = -30440
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3481.4546
$ws.Range("I22").Value = 4137
$ws.Range("K22").Value = 4137
$ws.Range("M22").Value = -3842
$ws.Range("H27").Value = 3481.4546
$ws.Range("I27").Value = 4137
$ws.Range("K27").Value = 4137
$ws.Range("M27").Value = -4030
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H93").Value = 403
$ws.Range("J93").Value = 454.5
$ws.Range("L93").Value = 454.5
$ws.Range("N93").Value = -2950.5
$ws.Range("H122").Value = 6332.2144
$ws.Range("I122").Value = 4856.5454
$ws.Range("J122").Value = 7287.0586
$ws.Range("K122").Value = 14569.6362
$ws.Range("L122").Value = 21861.1758
$ws.Range("M122").Value = -12119.6362
$ws.Range("N122").Value = -26761.1758
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9895.571
$ws.Range("J122").Value = 9924.25
$ws.Range("L122").Value = 29772.75
$ws.Range("N122").Value = -34672.75
$ws.Range("H132").Value = 3873.75
$ws.Range("I132").Value = 3873.75
$ws.Range("K132").Value = 11621.25
$ws.Range("M132").Value = -9091.25
$ws.Range("H136").Value = 1981.4
$ws.Range("I136").Value = 1369
$ws.Range("K136").Value = 4107
$ws.Range("M136").Value = -1557
